{"js": "// Replace each two-digit multiplication equation in the document with its\n// updated counterpart, per the commit diff. Every source string is unique\n// within the document, so a straightforward search+replace per pair is safe\n// and keeps the original run formatting (font, size, etc.) untouched since\n// insertText(\"Replace\") only rewrites the text inside the matched range.\nconst replacements = [\n  [\"99\u00d763=\", \"12\u00d761=\"],\n  [\"69\u00d753=\", \"48\u00d756=\"],\n  [\"85\u00d776=\", \"43\u00d758=\"],\n  [\"14\u00d784=\", \"92\u00d739=\"],\n  [\"81\u00d764=\", \"63\u00d761=\"],\n  [\"42\u00d742=\", \"23\u00d744=\"],\n  [\"54\u00d795=\", \"82\u00d714=\"],\n  [\"63\u00d740=\", \"23\u00d791=\"],\n  [\"85\u00d737=\", \"63\u00d795=\"],\n  [\"56\u00d754=\", \"27\u00d726=\"],\n  [\"54\u00d749=\", \"55\u00d783=\"],\n  [\"17\u00d734=\", \"93\u00d713=\"],\n  [\"73\u00d755=\", \"71\u00d759=\"],\n  [\"30\u00d750=\", \"37\u00d739=\"],\n  [\"60\u00d792=\", \"59\u00d712=\"],\n  [\"72\u00d799=\", \"38\u00d738=\"],\n  [\"19\u00d761=\", \"64\u00d780=\"],\n  [\"45\u00d745=\", \"35\u00d725=\"],\n  [\"70\u00d720=\", \"16\u00d723=\"],\n  [\"49\u00d792=\", \"12\u00d762=\"],\n  [\"87\u00d779=\", \"26\u00d736=\"],\n  [\"44\u00d750=\", \"45\u00d795=\"],\n  [\"97\u00d773=\", \"71\u00d773=\"],\n  [\"92\u00d791=\", \"61\u00d717=\"],\n  [\"21\u00d780=\", \"75\u00d721=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication equation in the document with its\n# updated counterpart, per the commit diff. Every source string is unique\n# within the document, so Find/Replace (wdReplaceAll) per pair is safe and\n# preserves the original run formatting (font, size, etc.) since Find.Execute\n# only rewrites the matched text.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @{ Old = \"99\u00d763=\"; New = \"12\u00d761=\" },\n  @{ Old = \"69\u00d753=\"; New = \"48\u00d756=\" },\n  @{ Old = \"85\u00d776=\"; New = \"43\u00d758=\" },\n  @{ Old = \"14\u00d784=\"; New = \"92\u00d739=\" },\n  @{ Old = \"81\u00d764=\"; New = \"63\u00d761=\" },\n  @{ Old = \"42\u00d742=\"; New = \"23\u00d744=\" },\n  @{ Old = \"54\u00d795=\"; New = \"82\u00d714=\" },\n  @{ Old = \"63\u00d740=\"; New = \"23\u00d791=\" },\n  @{ Old = \"85\u00d737=\"; New = \"63\u00d795=\" },\n  @{ Old = \"56\u00d754=\"; New = \"27\u00d726=\" },\n  @{ Old = \"54\u00d749=\"; New = \"55\u00d783=\" },\n  @{ Old = \"17\u00d734=\"; New = \"93\u00d713=\" },\n  @{ Old = \"73\u00d755=\"; New = \"71\u00d759=\" },\n  @{ Old = \"30\u00d750=\"; New = \"37\u00d739=\" },\n  @{ Old = \"60\u00d792=\"; New = \"59\u00d712=\" },\n  @{ Old = \"72\u00d799=\"; New = \"38\u00d738=\" },\n  @{ Old = \"19\u00d761=\"; New = \"64\u00d780=\" },\n  @{ Old = \"45\u00d745=\"; New = \"35\u00d725=\" },\n  @{ Old = \"70\u00d720=\"; New = \"16\u00d723=\" },\n  @{ Old = \"49\u00d792=\"; New = \"12\u00d762=\" },\n  @{ Old = \"87\u00d779=\"; New = \"26\u00d736=\" },\n  @{ Old = \"44\u00d750=\"; New = \"45\u00d795=\" },\n  @{ Old = \"97\u00d773=\"; New = \"71\u00d773=\" },\n  @{ Old = \"92\u00d791=\"; New = \"61\u00d717=\" },\n  @{ Old = \"21\u00d780=\"; New = \"75\u00d721=\" }\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute(\n    $pair.Old,  # FindText\n    $false,     # MatchCase\n    $false,     # MatchWholeWord\n    $false,     # MatchWildcards\n    $false,     # MatchSoundsLike\n    $false,     # MatchAllWordForms\n    $true,      # Forward\n    1,          # Wrap (wdFindContinue)\n    $false,     # Format\n    $pair.New,  # ReplaceWith\n    2           # Replace (wdReplaceAll)\n  ) | Out-Null\n}\n"}
